$d = $word.ActiveDocument

# --- 1. Fix spacing: "Personal id(either ID of a coach or the administrator)"
#        should read "Personal id (either ID of a coach or the administrator)"
$findRng = $d.Content
$findRng.Find.Execute("Personal id(either ID of a coach or the administrator)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Personal id (either ID of a coach or the administrator)", 2)

# --- 2. Word moves the hidden "_GoBack" bookmark to mark the location of the
#        last edit. Re-create it right after the word "Status" that follows
#        "Password" in the System Users bullet list (where the edit point was).
$pwdRng = $d.Content
$pwdRng.Find.Execute("Password")
$pwdPara = $pwdRng.Paragraphs.First

$allParas = $d.Paragraphs
$pwdStart = $pwdPara.Range.Start
$statusPara = $null
for ($i = 1; $i -le $allParas.Count; $i++) {
    $p = $allParas.Item($i)
    if ($p.Range.Start -eq $pwdStart) {
        $statusPara = $allParas.Item($i + 1)
        break
    }
}

if ($statusPara -ne $null) {
    $goBackPos = $statusPara.Range.End - 1
    $goBackRange = $d.Range($goBackPos, $goBackPos)

    try {
        $old = $d.Bookmarks.Item("_GoBack")
        $old.Delete()
    } catch {
    }

    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
